{"js": "// Di\u00e1rio de bordo do curso de HTML \u2014 \"salvando o diario de bordo\"\n//\n// The original paragraphs below had their text split across multiple\n// <w:r> runs with <w:proofErr> spell-check markers interleaved (an\n// artifact of Word's spell-checker flagging \"html\", \"head\", \"bory\",\n// \"lang\", \"pt\", \"tags\", \"title\", \"hr\" as unknown words). The edit\n// collapses each of those paragraphs back down to a single plain run\n// (no proofErr), and appends two new paragraphs of text after the\n// \"<hr> ...\" line (replacing the trailing blank paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraphs whose runs need to be merged back into a single run,\n// keyed by their (0-based) position in the document.\nconst mergedText = {\n  3: \"Criando o primeiro documento html, com </DOCTYPE HTML></HTML>\",\n  4: \"Dividido em head e body.\",\n  5: \"<head> cabe\u00e7alho onde tem as configura\u00e7\u00f5es iniciais.\",\n  6: \"<bory> o corpo onde ficam as informa\u00e7\u00f5es que o site vai apresentar.\",\n  7: \"Colocar em a lang em pt-br.\",\n  9: \"As tags;\",\n  10: \"<title> diz o nome do site\",\n  13: \"<hr> cria uma linha horizontal.\",\n};\n\nfor (const [index, text] of Object.entries(mergedText)) {\n  items[index].getRange().insertText(text, \"Replace\");\n}\n\n// The document's last paragraph is empty; turn it into the \"<br>\" line\n// and add one more paragraph after it for the \"&lt;/&gt;\" explanation.\nconst lastParagraph = items[items.length - 1];\nlastParagraph.getRange().insertText(\"<br> quebra de linha.\", \"Replace\");\nlastParagraph.insertParagraph(\n  \"Para mostrar os s\u00edmbolos de < = &lt;, e o de > &gt;\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Di\u00e1rio de bordo do curso de HTML \u2014 \"salvando o diario de bordo\"\n#\n# The affected paragraphs had their text split across several <w:r> runs\n# with <w:proofErr> spell-check markers in between (Word flagging \"html\",\n# \"head\", \"bory\", \"lang\", \"pt\", \"tags\", \"title\", \"hr\" as unknown words).\n# Running Find/Replace across each paragraph (old text -> same text)\n# forces Word to re-knit the paragraph into a single clean run and drops\n# the proofErr markers. Afterwards we turn the trailing blank paragraph\n# into the new \"<br> ...\" line and append one more paragraph of text.\n\n$d = $word.ActiveDocument\n\n$mergedText = [ordered]@{\n    4  = \"Criando o primeiro documento html, com </DOCTYPE HTML></HTML>\"\n    5  = \"Dividido em head e body.\"\n    6  = \"<head> cabe\u00e7alho onde tem as configura\u00e7\u00f5es iniciais.\"\n    7  = \"<bory> o corpo onde ficam as informa\u00e7\u00f5es que o site vai apresentar.\"\n    8  = \"Colocar em a lang em pt-br.\"\n    10 = \"As tags;\"\n    11 = \"<title> diz o nome do site\"\n    14 = \"<hr> cria uma linha horizontal.\"\n}\n\nforeach ($index in $mergedText.Keys) {\n    $text = $mergedText[$index]\n    $find = $d.Paragraphs($index).Range.Find\n    $find.ClearFormatting()\n    $find.Text = $text\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $text\n    $find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null\n}\n\n# The last paragraph in the document is empty; fill it with the \"<br>\"\n# line, then add a new paragraph after it with the symbols explanation.\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$lastParagraph.Range.Text = \"<br> quebra de linha.\"\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$newParagraph.Range.Text = \"Para mostrar os s\u00edmbolos de < = &lt;, e o de > &gt;\"\n"}
